$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.0959030901961
$ws.Range("C2").Value = 6.087916234589967
$ws.Range("D2").Value = 13.77070329496899
$ws.Range("E2").Value = 14.18009475633176
$ws.Range("G2").Value = 3.733869478229223
$ws.Range("I2").Value = 33.02233627386346
$ws.Range("J2").Value = 8.586003250324362
$ws.Range("K2").Value = 11.47194820863597
$ws.Range("L2").Value = 12.55167725221384
$ws.Range("M2").Value = 16.60304737728586
$ws.Range("O2").Value = 34.52129130251089
$ws.Range("B3").Value = 13.949861615546
$ws.Range("C3").Value = 6.016886658555023
$ws.Range("D3").Value = 13.77405470356705
$ws.Range("E3").Value = 14.20584320711117
$ws.Range("G3").Value = 3.735927541434649
$ws.Range("I3").Value = 33.1125060136506
$ws.Range("J3").Value = 8.590414220800605
$ws.Range("K3").Value = 11.37005805231713
$ws.Range("L3").Value = 12.56557178713197
$ws.Range("M3").Value = 16.59265597018036
$ws.Range("O3").Value = 34.60510711114545
$ws.Range("B4").Value = 13.8620809018982
$ws.Range("C4").Value = 5.972122363195941
$ws.Range("D4").Value = 13.77835391514031
$ws.Range("E4").Value = 14.22311629173331
$ws.Range("G4").Value = 3.737258959061593
$ws.Range("I4").Value = 33.17258247111704
$ws.Range("J4").Value = 8.593273496178687
$ws.Range("K4").Value = 11.30905276538524
$ws.Range("L4").Value = 12.57556202440454
$ws.Range("M4").Value = 16.58848715258694
$ws.Range("O4").Value = 34.66196454926108
$ws.Range("B5").Value = 13.82682250521217
$ws.Range("C5").Value = 5.953599063947452
$ws.Range("D5").Value = 13.78067052959611
$ws.Range("E5").Value = 14.23052372241373
$ws.Range("G5").Value = 3.737818614547472
$ws.Range("I5").Value = 33.19824889003227
$ws.Range("J5").Value = 8.594476731202944
$ws.Range("K5").Value = 11.28460793836995
$ws.Range("L5").Value = 12.58000054918493
$ws.Range("M5").Value = 16.5873471434315
$ws.Range("O5").Value = 34.68649005413064
$ws.Range("B6").Value = 13.82099994512904
$ws.Range("C6").Value = 5.950506525300685
$ws.Range("D6").Value = 13.7810893313031
$ws.Range("E6").Value = 14.23177599265868
$ws.Range("G6").Value = 3.737912578805403
$ws.Range("I6").Value = 33.20258233476777
$ws.Range("J6").Value = 8.594678829156344
$ws.Range("K6").Value = 11.28057466894442
$ws.Range("L6").Value = 12.58075976986848
$ws.Range("M6").Value = 16.58719167171315
$ws.Range("O6").Value = 34.69064434858588
$ws.Range("B7").Value = 13.86160326854123
$ws.Range("C7").Value = 5.97187368030993
$ws.Range("D7").Value = 13.77838287034904
$ws.Range("E7").Value = 14.22321469817336
$ws.Range("G7").Value = 3.737266437495482
$ws.Range("I7").Value = 33.17292382042217
$ws.Range("J7").Value = 8.593289569192535
$ws.Range("K7").Value = 11.30872138092924
$ws.Range("L7").Value = 12.57562039550971
$ws.Range("M7").Value = 16.58846951198312
$ws.Range("O7").Value = 34.6622898209759
$ws.Range("B8").Value = 14.04517883775925
$ws.Range("C8").Value = 6.063669129275216
$ws.Range("D8").Value = 13.77139432019922
$ws.Range("E8").Value = 14.18866939935208
$ws.Range("G8").Value = 3.734565067438838
$ws.Range("I8").Value = 33.05244881861941
$ws.Range("J8").Value = 8.587492893403166
$ws.Range("K8").Value = 11.43650829094896
$ws.Range("L8").Value = 12.55616563190797
$ws.Range("M8").Value = 16.5990070985289
$ws.Range("O8").Value = 34.54907098318182
$ws.Range("B9").Value = 14.41846000880712
$ws.Range("C9").Value = 6.234194714219348
$ws.Range("D9").Value = 13.77542295624769
$ws.Range("E9").Value = 14.13251494446759
$ws.Range("G9").Value = 3.729802879202341
$ws.Range("I9").Value = 32.85358688647961
$ws.Range("J9").Value = 8.577318268735159
$ws.Range("K9").Value = 11.69835078823984
$ws.Range("L9").Value = 12.52956558441183
$ws.Range("M9").Value = 16.63709261439818
$ws.Range("O9").Value = 34.36988279828265
$ws.Range("B10").Value = 14.69849800801308
$ws.Range("C10").Value = 6.35325369567582
$ws.Range("D10").Value = 13.78911432428225
$ws.Range("E10").Value = 14.09829197607255
$ws.Range("G10").Value = 3.726626947940606
$ws.Range("I10").Value = 32.73027952376165
$ws.Range("J10").Value = 8.570563249928433
$ws.Range("K10").Value = 11.89609169252858
$ws.Range("L10").Value = 12.51702902777796
$ws.Range("M10").Value = 16.67551426252225
$ws.Range("O10").Value = 34.26438852408577
$ws.Range("B11").Value = 14.82665608484369
$ws.Range("C11").Value = 6.405975652978574
$ws.Range("D11").Value = 13.79765211309132
$ws.Range("E11").Value = 14.08424362643913
$ws.Range("G11").Value = 3.725251511397556
$ws.Range("I11").Value = 32.67913437328739
$ws.Range("J11").Value = 8.567645178438758
$ws.Range("K11").Value = 11.98689018003495
$ws.Range("L11").Value = 12.51283894114344
$ws.Range("M11").Value = 16.69521726724652
$ws.Range("O11").Value = 34.22208347690047
$ws.Range("B12").Value = 14.87525530570956
$ws.Range("C12").Value = 6.425726020775048
$ws.Range("D12").Value = 13.80121508792275
$ws.Range("E12").Value = 14.07914187855223
$ws.Range("G12").Value = 3.724740581062881
$ws.Range("I12").Value = 32.66047863706821
$ws.Range("J12").Value = 8.566562335003264
$ws.Range("K12").Value = 12.02136742029214
$ws.Range("L12").Value = 12.51146902458705
$ws.Range("M12").Value = 16.70299414599928
$ws.Range("O12").Value = 34.20688162200916
$ws.Range("B13").Value = 14.8647862057851
$ws.Range("C13").Value = 6.421482071655022
$ws.Range("D13").Value = 13.80043310261073
$ws.Range("E13").Value = 14.08023094108179
$ws.Range("G13").Value = 3.724850178741216
$ws.Range("I13").Value = 32.66446482049934
$ws.Range("J13").Value = 8.566794560352982
$ws.Range("K13").Value = 12.01393839666677
$ws.Range("L13").Value = 12.51175443210946
$ws.Range("M13").Value = 16.70130528011269
$ws.Range("O13").Value = 34.21011921640368
$ws.Range("B14").Value = 14.83065323601208
$ws.Range("C14").Value = 6.40760486390527
$ws.Range("D14").Value = 13.7979386429043
$ws.Range("E14").Value = 14.08381953525067
$ws.Range("G14").Value = 3.725209278367125
$ws.Range("I14").Value = 32.67758528697849
$ws.Range("J14").Value = 8.567555648547156
$ws.Range("K14").Value = 11.98972492122933
$ws.Range("L14").Value = 12.51272189729804
$ws.Range("M14").Value = 16.69585077162605
$ws.Range("O14").Value = 34.2208164107284
$ws.Range("B15").Value = 14.80975350006
$ws.Range("C15").Value = 6.399076550598972
$ws.Range("D15").Value = 13.79645360778612
$ws.Range("E15").Value = 14.08604603541701
$ws.Range("G15").Value = 3.725430527373693
$ws.Range("I15").Value = 32.68571466057465
$ws.Range("J15").Value = 8.568024720991298
$ws.Range("K15").Value = 11.97490486672288
$ws.Range("L15").Value = 12.51334270496041
$ws.Range("M15").Value = 16.69255072443403
$ws.Range("O15").Value = 34.22747531652528
$ws.Range("B16").Value = 14.69013443296978
$ws.Range("C16").Value = 6.349778602519009
$ws.Range("D16").Value = 13.78860264945105
$ws.Range("E16").Value = 14.09924060990722
$ws.Range("G16").Value = 3.726718226448623
$ws.Range("I16").Value = 32.73372157257439
$ws.Range("J16").Value = 8.570757059911777
$ws.Range("K16").Value = 11.89017243354663
$ws.Range("L16").Value = 12.51733323773258
$ws.Range("M16").Value = 16.67427104441673
$ws.Range("O16").Value = 34.26726771188881
$ws.Range("B17").Value = 14.61691866246633
$ws.Range("C17").Value = 6.31916207992132
$ws.Range("D17").Value = 13.7843763626821
$ws.Range("E17").Value = 14.10772397468658
$ws.Range("G17").Value = 3.727525905121121
$ws.Range("I17").Value = 32.76443978209633
$ws.Range("J17").Value = 8.572472845914962
$ws.Range("K17").Value = 11.83838855260876
$ws.Range("L17").Value = 12.52016824756421
$ws.Range("M17").Value = 16.66362400173009
$ws.Range("O17").Value = 34.29313559238527
$ws.Range("B18").Value = 14.57488159617292
$ws.Range("C18").Value = 6.301417287309715
$ws.Range("D18").Value = 13.78216304277699
$ws.Range("E18").Value = 14.11274646747119
$ws.Range("G18").Value = 3.727996987276864
$ws.Range("I18").Value = 32.78257379455141
$ws.Range("J18").Value = 8.573474298764641
$ws.Range("K18").Value = 11.80868504022804
$ws.Range("L18").Value = 12.52194131640812
$ws.Range("M18").Value = 16.65770977617075
$ws.Range("O18").Value = 34.30854916170446
$ws.Range("B19").Value = 14.56066268071117
$ws.Range("C19").Value = 6.295386253740063
$ws.Range("D19").Value = 13.781451075073
$ws.Range("E19").Value = 14.11447158727911
$ws.Range("G19").Value = 3.728157610079336
$ws.Range("I19").Value = 32.78879363616917
$ws.Range("J19").Value = 8.573815880241062
$ws.Range("K19").Value = 11.79864273443341
$ws.Range("L19").Value = 12.52256613494473
$ws.Range("M19").Value = 16.65574345491568
$ws.Range("O19").Value = 34.31385980428364
$ws.Range("B20").Value = 14.62470517860718
$ws.Range("C20").Value = 6.322435285518592
$ws.Range("D20").Value = 13.78480376004294
$ws.Range("E20").Value = 14.10680610074738
$ws.Range("G20").Value = 3.727439251211202
$ws.Range("I20").Value = 32.7611215746134
$ws.Range("J20").Value = 8.572288689484129
$ws.Range("K20").Value = 11.84389283165236
$ws.Range("L20").Value = 12.51985171869799
$ws.Range("M20").Value = 16.66473572666136
$ws.Range("O20").Value = 34.29032653156569
$ws.Range("B21").Value = 14.84067738713768
$ws.Range("C21").Value = 6.41168680984091
$ws.Range("D21").Value = 13.79866239146564
$ws.Range("E21").Value = 14.08275956414034
$ws.Range("G21").Value = 3.725103533329899
$ws.Range("I21").Value = 32.67371216693123
$ws.Range("J21").Value = 8.567331497653271
$ws.Range("K21").Value = 11.99683467974659
$ws.Range("L21").Value = 12.5124318523843
$ws.Range("M21").Value = 16.69744435701791
$ws.Range("O21").Value = 34.2176521762548
$ws.Range("B22").Value = 14.98220822212742
$ws.Range("C22").Value = 6.468765534607092
$ws.Range("D22").Value = 13.80964153890753
$ws.Range("E22").Value = 14.06831452998628
$ws.Range("G22").Value = 3.723634791712636
$ws.Range("I22").Value = 32.62073423807321
$ws.Range("J22").Value = 8.564220847976271
$ws.Range("K22").Value = 12.09732516090292
$ws.Range("L22").Value = 12.50884569982284
$ws.Range("M22").Value = 16.72065985369347
$ws.Range("O22").Value = 34.17492431898373
$ws.Range("B23").Value = 14.90664940363192
$ws.Range("C23").Value = 6.438418510390568
$ws.Range("D23").Value = 13.80360670757475
$ws.Range("E23").Value = 14.0759080074508
$ws.Range("G23").Value = 3.724413415550776
$ws.Range("I23").Value = 32.6486298102872
$ws.Range("J23").Value = 8.56586927345416
$ws.Range("K23").Value = 12.04365174026906
$ws.Range("L23").Value = 12.51064438426475
$ws.Range("M23").Value = 16.70810251003322
$ws.Range("O23").Value = 34.19729241670669
$ws.Range("B24").Value = 14.6211847181067
$ws.Range("C24").Value = 6.320955913969792
$ws.Range("D24").Value = 13.78460985938922
$ws.Range("E24").Value = 14.10722061894759
$ws.Range("G24").Value = 3.727478406455344
$ws.Range("I24").Value = 32.76262026064995
$ws.Range("J24").Value = 8.572371899799357
$ws.Range("K24").Value = 11.84140413400839
$ws.Range("L24").Value = 12.51999437533826
$ws.Range("M24").Value = 16.66423247103198
$ws.Range("O24").Value = 34.29159482037045
$ws.Range("B25").Value = 14.31629254942131
$ws.Range("C25").Value = 6.189125214296443
$ws.Range("D25").Value = 13.77244197614271
$ws.Range("E25").Value = 14.14646873810921
$ws.Range("G25").Value = 3.731034236293076
$ws.Range("I25").Value = 32.90338116659447
$ws.Range("J25").Value = 8.5799437989795
$ws.Range("K25").Value = 11.62646367034998
$ws.Range("L25").Value = 12.53552847117025
$ws.Range("M25").Value = 16.62494208533576
$ws.Range("O25").Value = 34.41376756072101
